$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace sample data with "optional" / "requied" placeholders ---
$ws.Range("A2").Value = "optional"
$ws.Range("B2").Value = "requied"
$ws.Range("C2").Value = "requied"
$ws.Range("D2").Value = "optional"
$ws.Range("E2").Value = "optional"
$ws.Range("F2").Value = "requied"
$ws.Range("G2").Value = "requied"
$ws.Range("H2").Value = "requied"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = "optional"
$ws.Range("K2").Value = "optional"
$ws.Range("L2").Value = "optional"
$ws.Range("M2").Value = "optional"

# --- Row 3: new book record (order matches the original authoring so the
# regenerated shared-string table lines up with the source workbook) ---
$ws.Range("A3").Value = 1000
$ws.Range("B3").Value = "learning C"
$ws.Range("C3").Value = "denesh Pathak "
$ws.Range("D3").Value = "k ramanujan"
$ws.Range("G3").Value = "ph pub."
$ws.Range("L3").Value = "main library"
$ws.Range("E3").ClearContents()
$ws.Range("K3").ClearContents()
$ws.Range("N3").Value = 40

# --- Row 4: another book record ---
$ws.Range("B4").Value = "Java"
$ws.Range("C4").Value = "zakie m"
$ws.Range("F4").Value = "7th"
$ws.Range("F3").Value = "2nd"
$ws.Range("G4").Value = "jk publishers"
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()

# --- Row 5 no longer exists in the rolled-back sheet ---
$ws.Rows("5:5").Delete()

# --- Column widths ---
$ws.Columns("B:B").ColumnWidth = 28.830729166666668
$ws.Columns("C:C").ColumnWidth = 25.608072916666668
$ws.Columns("D:D").ColumnWidth = 32.276041666666664
$ws.Columns("E:E").ColumnWidth = 24.276041666666668
$ws.Columns("F:F").ColumnWidth = 13.053385416666666
$ws.Columns("G:G").ColumnWidth = 22.276041666666668
$ws.Columns("H:H").ColumnWidth = 12.166666666666666
$ws.Columns("I:I").ColumnWidth = 11.498697916666666
$ws.Columns("K:M").ColumnWidth = 12.830729166666666

# --- Selection moves to D13 ---
$ws.Range("D13").Select()
